# Updates crypto price/volume data per the Wed May 29 17:52:43 UTC 2024 GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.360.20'
$ws.Range("E2").Value = '  -0.97%  '
$ws.Range("D3").Value = '3.754.96'
$ws.Range("E3").Value = '  -2.10%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = "'594.80"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.94%  '
$ws.Range("D6").Value = "'169.32"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.16%  '
$ws.Range("D7").Value = '3.753.11'
$ws.Range("E7").Value = '  -2.19%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("D9").Value = "'0.527"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.30%  '
$ws.Range("D10").Value = "'0.165"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.44%  '
$ws.Range("D11").Value = "'6.47"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.06%  '
$ws.Range("D12").Value = "'0.454"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.98%  '
$ws.Range("D13").Value = "'0.0000276"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +6.09%  '
$ws.Range("D14").Value = "'36.58"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.46%  '
$ws.Range("D15").Value = '4.400.57'
$ws.Range("E15").Value = '  -1.91%  '
$ws.Range("D16").Value = '3.764.14'
$ws.Range("E16").Value = '  -2.03%  '
$ws.Range("D17").Value = "'18.75"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.47%  '
$ws.Range("D18").Value = '67.502.45'
$ws.Range("E18").Value = '  -0.96%  '
$ws.Range("D19").Value = "'7.21"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.60%  '
$ws.Range("E20").Value = '  +0.86%  '
$ws.Range("D21").Value = "'10.54"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -5.29%  '
$ws.Range("D22").Value = "'468.41"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.45%  '
$ws.Range("D23").Value = "'0.719"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.12%  '
$ws.Range("D24").Value = "'0.0000147"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -8.04%  '
$ws.Range("D25").Value = "'83.79"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.10%  '
$ws.Range("D26").Value = "'2.22"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.36%  '
$ws.Range("D27").Value = "'12.12"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.14%  '
$ws.Range("D28").Value = "'10.36"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.15%  '
$ws.Range("D30").Value = "'2.90"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.23%  '
$ws.Range("D31").Value = '3.916.15'
$ws.Range("E31").Value = '  -1.86%  '
$ws.Range("D32").Value = "'7.69"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.51%  '
$ws.Range("B33").Value = 'EthereumClassic'
$ws.Range("C33").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D33").Value = "'30.39"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.81%  '
$ws.Range("B34").Value = 'ImmutableX'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D34").Value = "'2.24"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.41%  '
$ws.Range("D35").Value = "'9.12"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.82%  '
$ws.Range("D36").Value = '3.732.10'
$ws.Range("D37").Value = "'3.82"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +4.44%  '
$ws.Range("D38").Value = "'0.105"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.59%  '
$ws.Range("E39").Value = '  -1.84%  '
$ws.Range("D40").Value = "'5.86"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.06%  '
$ws.Range("D41").Value = "'0.998"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.90%  '
$ws.Range("E42").Value = '  +0.10%  '
$ws.Range("D43").Value = "'0.312"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.64%  '
$ws.Range("D45").Value = "'8.73"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.22%  '
$ws.Range("D46").Value = "'1.95"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.71%  '
$ws.Range("D47").Value = "'45.81"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.36%  '
$ws.Range("D48").Value = "'396.12"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -5.92%  '
$ws.Range("B49").Value = 'Monero'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D49").Value = "'141.69"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.15%  '
$ws.Range("B50").Value = 'FLOKI'
$ws.Range("C50").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D50").Value = "'0.000268"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -8.61%  '
$ws.Range("D51").Value = "'0.0354"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.93%  '
